# LevelDbSheet.xlsx update
# - Adds a new "enemyType" column (H) classifying each spawn row as Minion/Boss
# - Normalizes spawnCount / spawnDelayTime / countDownTime values
# - Updates the selected cell to E14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H
$ws.Range("H1").Value = "enemyType"

# Data rows: dev_ID, spawnCount, spawnDelayTime, countDownTime, spawnType1, spawnType2, spawnType3, enemyType
$data = @(
    @(1,  10, 2, 10, "Minion_Ball",   "None",          "None",          "Minion"),
    @(2,  10, 2, 10, "Minion_Spider", "None",          "None",          "Minion"),
    @(3,  10, 2, 10, "Minion_Drone",  "None",          "None",          "Minion"),
    @(4,  10, 2, 10, "Minion_Ball",   "Minion_Spider", "None",          "Minion"),
    @(5,  10, 2, 10, "Minion_Ball",   "Minion_Spider", "None",          "Minion"),
    @(6,  10, 2, 10, "Minion_Ball",   "Minion_Drone",  "None",          "Minion"),
    @(7,  10, 2, 10, "Minion_Ball",   "Minion_Drone",  "None",          "Minion"),
    @(8,  20, 2, 10, "Minion_Ball",   "Minion_Spider", "Minion_Drone",  "Minion"),
    @(9,  20, 2, 10, "Minion_Ball",   "Minion_Spider", "Minion_Drone",  "Minion"),
    @(10, 1,  0, 120, "Boss_SkyFire", "None",          "None",          "Boss")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
    $ws.Cells.Item($row, 8).Value = $values[7]
}

# New column width for H
$ws.Columns.Item(8).ColumnWidth = 11.5

# Update selection to match the authored edit
$ws.Range("E14").Select()
